$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 2 (bold "Required setup" style row, copied from row 1 above) ---
$ws.Rows(2).Insert()
$ws.Range("A2").Value = "n/a"
$ws.Range("B2").Value = "addi `$t2, `$zero, 7"
$ws.Range("C2:D2").Clear()

# --- Append a new row 8 with the div test case ---
$ws.Range("B8").Value = "div `$t2, `$t0"
$ws.Range("A8").Value = "Set `$t2 to 7; set `$t0 to 6"
$ws.Range("D8").Value = "# LO = 7 // 6 = 1, HI = 7 % 6 = 1"
$ws.Range("C8").Value = "0x0148001a"

# --- View state: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 165
$ws.Range("C13").Select()
